# Commit: "Fruta / hortaliza, semanal"
# A new weekly price record is inserted into the price history table.
# This shifts all existing data rows from row 10 onward down by one row,
# and the new row 10 receives the latest week's data (Fecha = 44764,
# Volumen = 50, everything else identical to what used to be in row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above row 10 - this pushes the previous
# row 10..68 down to 11..69, automatically carrying over cell styles
# (e.g. the date format on column D) and extending the sheet dimension.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C10").Value = "Los Lagos"
$ws.Range("D10").Value = 44764
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 100112043
$ws.Range("G10").Value = "Pepino dulce"
$ws.Range("H10").Value = "Cultivar IV Región"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 18000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 18000
$ws.Range("N10").Value = "`$/bandeja 18 kilos"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 1000
$ws.Range("Q10").Value = 18
$ws.Range("R10").Value = "Hortaliza"
